# Update the "last_edited_time" (column D) values for rows 2-83 in the
# LUY_KE_NGAY_LONG_XUYEN sheet. The original timestamps on 2024-07-18
# (15:58 / 15:59 / 16:00) are replaced with new timestamps on 2024-07-19
# (08:00 / 08:01 / 08:02 / 08:03), preserving the minute-by-minute
# progression as the row number increases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Start = 2;  End = 3;  Value = "2024-07-19T08:00:00.000Z" },
    @{ Start = 4;  End = 25; Value = "2024-07-19T08:01:00.000Z" },
    @{ Start = 26; End = 58; Value = "2024-07-19T08:02:00.000Z" },
    @{ Start = 59; End = 83; Value = "2024-07-19T08:03:00.000Z" }
)

foreach ($u in $updates) {
    for ($r = $u.Start; $r -le $u.End; $r++) {
        $ws.Cells.Item($r, 4).Value = $u.Value
    }
}
